$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" (column F) values for each row, per repulled data.
$updates = @{
    2  = -9
    3  = -2
    4  = -2
    6  = -1
    7  = 6
    8  = 1
    9  = -2
    10 = -2
    11 = -1
    12 = -5
    13 = -2
    14 = -2
    15 = 2
    16 = -2
    17 = 3
    18 = 1
    19 = -3
    21 = -2
    22 = 5
    23 = 1
    24 = 8
    25 = -1
    26 = -6
    27 = 3
    28 = -1
    29 = 12
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
